# Bids_template.xlsx change:
#   "Show user full name rather than login name for all kinds of data."
#
# The template's "sales person" column currently renders the JXLS
# expression ${record.salesPerson} (login name). Replace it with
# ${record.salesPersonFullName} so the generated report shows the
# salesperson's full name instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 holds the per-record template placeholder for the "sales person"
# column (row 1 is the header, row 2 is the JXLS "each" template row).
$ws.Range("D2").Value2 = '${record.salesPersonFullName}'

# Reflect the cell the author was last working on/reviewing.
$ws.Range("D3").Select()
